# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with refreshed values from the data source.
#
# The Price/Volume cells are stored as *text* (many prices contain
# thousands separators written as extra dots, e.g. "29.303.26", which are
# not valid numbers), so values that happen to look numeric are written
# with a leading apostrophe to force Excel to keep them as text, then the
# cell style is reset to "Normal" so no stray NumberFormat/quote-prefix
# style is left attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '29.303.26'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '1.845.12'
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("D4").Value = '''0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''240.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.23%  '
$ws.Range("D6").Value = '''0.6727'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.96%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '''0.07449'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '''0.2949'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.39%  '
$ws.Range("D10").Value = '''22.97'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("D12").Value = '1.844.26'
$ws.Range("E12").Value = '  +0.11%  '
$ws.Range("D13").Value = '''5.010'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.00%  '
$ws.Range("D14").Value = '''0.6721'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = '''85.99'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.82%  '
$ws.Range("D16").Value = '''6.152'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.37%  '
$ws.Range("D17").Value = '29.285.11'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("D18").Value = '''0.000008325'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.99%  '
$ws.Range("D19").Value = '''229.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.67%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").Value = '''7.189'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.73%  '
$ws.Range("D23").Value = '''1.000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").Value = '''161.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.53%  '
$ws.Range("D25").Value = '''8.714'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.49%  '
$ws.Range("D26").Value = '''0.1409'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("D27").Value = '''18.03'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.39%  '
$ws.Range("D28").Value = '''1.512'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("D29").Value = '''4.167'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.34%  '
$ws.Range("D30").Value = '''4.072'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.60%  '
$ws.Range("D31").Value = '''1.194'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.04%  '
$ws.Range("D32").Value = '''0.05311'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.39%  '
$ws.Range("D33").Value = '''0.7581'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.62%  '
$ws.Range("D34").Value = '''1.876'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.40%  '
$ws.Range("D35").Value = '''1.139'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("D37").Value = '1.321.97'
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("D38").Value = '''0.01805'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("D39").Value = '''2.725'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").Value = '''0.9207'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("D41").Value = '''5.987'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.24%  '
$ws.Range("D42").Value = '''0.08389'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +12.92%  '
$ws.Range("D43").Value = '''1.002'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("D44").Value = '''103.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.66%  '
$ws.Range("D45").Value = '1.991.25'
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("D46").Value = '''0.00000000123'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("D47").Value = '''0.5167'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.70%  '
$ws.Range("D48").Value = '''1.779'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("D49").Value = '''64.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.34%  '
$ws.Range("D50").Value = '''9.127'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.16%  '
$ws.Range("D51").Value = '''0.05954'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.17%  '
